$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 946.2727
$ws.Range("I32").Value = 677.1429000000001
$ws.Range("J32").Value = 1071.8667
$ws.Range("K32").Value = 677.1429000000001
$ws.Range("L32").Value = 1071.8667
$ws.Range("M32").Value = -351.1429000000001
$ws.Range("N32").Value = -1723.8667
$ws.Range("H106").Value = 121213710
$ws.Range("I106").Value = 33335080
$ws.Range("J106").Value = 1000000000
$ws.Range("K106").Value = 33335080
$ws.Range("L106").Value = 1000000000
$ws.Range("M106").Value = -33334449
$ws.Range("H137").Value = 2600
$ws.Range("I137").Value = 2747.353
$ws.Range("J137").Value = 1973.75
$ws.Range("K137").Value = 8242.059000000001
$ws.Range("L137").Value = 5921.25
$ws.Range("M137").Value = -5692.059000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6954.25
$ws.Range("I61").Value = 8339
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 8339
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -8127
$ws.Range("N61").Value = -3224
$ws.Range("H74").Value = 1742.5
$ws.Range("I74").Value = 1631.0613
$ws.Range("J74").Value = 2238.9092
$ws.Range("K74").Value = 1631.0613
$ws.Range("L74").Value = 2238.9092
$ws.Range("M74").Value = -757.0613000000001
$ws.Range("N74").Value = -3986.9092
$ws.Range("H77").Value = 1742.5
$ws.Range("I77").Value = 1631.0613
$ws.Range("J77").Value = 2238.9092
$ws.Range("K77").Value = 8155.306500000001
$ws.Range("L77").Value = 11194.546
$ws.Range("M77").Value = -3787.306500000001
$ws.Range("N77").Value = -19930.546
$ws.Range("H97").Value = 808.7143
$ws.Range("I97").Value = 608.3333
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 608.3333
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -112.3333
$ws.Range("N97").Value = -3003
$ws.Range("H106").Value = 59444
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 59444
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 59444
$ws.Range("N106").Value = -61968
$ws.Range("H109").Value = 49000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 49000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 49000
$ws.Range("N109").Value = -51774
$ws.Range("H136").Value = 6954.25
$ws.Range("I136").Value = 8339
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 25017
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -22467
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 30499
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 30499
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 30499
$ws.Range("N51").Value = -31481
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H107").Value = 1150.4117
$ws.Range("I107").Value = 1049.6
$ws.Range("J107").Value = 1906.5
$ws.Range("K107").Value = 1049.6
$ws.Range("L107").Value = 1906.5
$ws.Range("M107").Value = 870.4000000000001
$ws.Range("N107").Value = -5746.5
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H134").Value = 3979.5957
$ws.Range("I134").Value = 4303.4326
$ws.Range("J134").Value = 2781.4
$ws.Range("K134").Value = 12910.2978
$ws.Range("L134").Value = 8344.200000000001
$ws.Range("M134").Value = -10375.2978
$ws.Range("N134").Value = -13414.2
$ws.Range("H138").Value = 53550
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 53550
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 53550
$ws.Range("N138").Value = -63830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10924.143
$ws.Range("I31").Value = 1008
$ws.Range("J31").Value = 13628.546
$ws.Range("K31").Value = 1008
$ws.Range("L31").Value = 13628.546
$ws.Range("M31").Value = -713
$ws.Range("N31").Value = -14218.546
$ws.Range("H34").Value = 10924.143
$ws.Range("I34").Value = 1008
$ws.Range("J34").Value = 13628.546
$ws.Range("K34").Value = 1008
$ws.Range("L34").Value = 13628.546
$ws.Range("M34").Value = -806
$ws.Range("N34").Value = -14032.546
$ws.Range("H58").Value = 962.0213
$ws.Range("I58").Value = 797.2
$ws.Range("J58").Value = 1442.75
$ws.Range("K58").Value = 797.2
$ws.Range("L58").Value = 1442.75
$ws.Range("M58").Value = -594.2
$ws.Range("N58").Value = -1848.75
$ws.Range("H99").Value = 3293.2354
$ws.Range("I99").Value = 1477.2222
$ws.Range("J99").Value = 5336.25
$ws.Range("K99").Value = 1477.2222
$ws.Range("L99").Value = 5336.25
$ws.Range("M99").Value = 20.77780000000007
$ws.Range("N99").Value = -8332.25
$ws.Range("H109").Value = 50190
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 50190
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 50190
$ws.Range("N109").Value = -52270
$ws.Range("H126").Value = 3293.2354
$ws.Range("I126").Value = 1477.2222
$ws.Range("J126").Value = 5336.25
$ws.Range("K126").Value = 4431.6666
$ws.Range("L126").Value = 16008.75
$ws.Range("M126").Value = -1961.6666
$ws.Range("N126").Value = -20948.75
$ws.Range("H132").Value = 1886.8572
$ws.Range("I132").Value = 1378
$ws.Range("J132").Value = 4049.5
$ws.Range("K132").Value = 4134
$ws.Range("L132").Value = 12148.5
$ws.Range("M132").Value = -1604
$ws.Range("N132").Value = -17208.5
$ws.Range("H134").Value = 1439.8422
$ws.Range("I134").Value = 1373.9412
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 4121.8236
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -1586.8236
$ws.Range("H136").Value = 962.0213
$ws.Range("I136").Value = 797.2
$ws.Range("J136").Value = 1442.75
$ws.Range("K136").Value = 2391.6
$ws.Range("L136").Value = 4328.25
$ws.Range("M136").Value = 158.3999999999996
$ws.Range("N136").Value = -9428.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 334242.78
$ws.Range("I5").Value = 604.75
$ws.Range("J5").Value = 601153.2
$ws.Range("K5").Value = 1814.25
$ws.Range("L5").Value = 1803459.6
$ws.Range("M5").Value = -1702.25
$ws.Range("N5").Value = -1803683.6
$ws.Range("H109").Value = 2085.6667
$ws.Range("I109").Value = 860.7143
$ws.Range("J109").Value = 3157.5
$ws.Range("K109").Value = 2582.1429
$ws.Range("L109").Value = 9472.5
$ws.Range("M109").Value = -1542.1429
$ws.Range("N109").Value = -11552.5
$ws.Range("H112").Value = 2582.4814
$ws.Range("I112").Value = 1027
$ws.Range("J112").Value = 2642.3076
$ws.Range("K112").Value = 3081
$ws.Range("L112").Value = 7926.9228
$ws.Range("M112").Value = -1973
$ws.Range("N112").Value = -10142.9228
$ws.Range("H122").Value = 742.0714
$ws.Range("I122").Value = 683.8570999999999
$ws.Range("J122").Value = 800.2857
$ws.Range("K122").Value = 6154.7139
$ws.Range("L122").Value = 7202.571300000001
$ws.Range("M122").Value = -3704.7139
$ws.Range("N122").Value = -12102.5713
$ws.Range("H135").Value = 334242.78
$ws.Range("I135").Value = 604.75
$ws.Range("J135").Value = 601153.2
$ws.Range("K135").Value = 5442.75
$ws.Range("L135").Value = 5410378.8
$ws.Range("M135").Value = -2907.75
$ws.Range("N135").Value = -5415448.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 10292.75
$ws.Range("I126").Value = 12879.111
$ws.Range("J126").Value = 2533.6667
$ws.Range("K126").Value = 38637.333
$ws.Range("L126").Value = 7601.000100000001
$ws.Range("M126").Value = -36167.333
$ws.Range("N126").Value = -12541.0001
$ws.Range("H132").Value = 2678.1409
$ws.Range("I132").Value = 2400.425
$ws.Range("J132").Value = 3036.484
$ws.Range("K132").Value = 7201.275000000001
$ws.Range("L132").Value = 9109.451999999999
$ws.Range("M132").Value = -4671.275000000001
$ws.Range("N132").Value = -14169.452

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3586738.2
$ws.Range("I22").Value = 22222662
$ws.Range("J22").Value = 2906.7693
$ws.Range("K22").Value = 22222662
$ws.Range("L22").Value = 2906.7693
$ws.Range("M22").Value = -22222367
$ws.Range("N22").Value = -3496.7693
$ws.Range("H27").Value = 3586738.2
$ws.Range("I27").Value = 22222662
$ws.Range("J27").Value = 2906.7693
$ws.Range("K27").Value = 22222662
$ws.Range("L27").Value = 2906.7693
$ws.Range("M27").Value = -22222555
$ws.Range("N27").Value = -3120.7693
$ws.Range("H46").Value = 30303832
$ws.Range("I46").Value = 55556132
$ws.Range("J46").Value = 1070
$ws.Range("K46").Value = 55556132
$ws.Range("L46").Value = 1070
$ws.Range("M46").Value = -55555944
$ws.Range("H132").Value = 13547417
$ws.Range("I132").Value = 17339654
$ws.Range("J132").Value = 3714.7144
$ws.Range("K132").Value = 52018962
$ws.Range("L132").Value = 11144.1432
$ws.Range("M132").Value = -52016432
$ws.Range("N132").Value = -16204.1432
$ws.Range("H136").Value = 11984.083
$ws.Range("I136").Value = 21149.834
$ws.Range("J136").Value = 2818.3333
$ws.Range("K136").Value = 63449.50199999999
$ws.Range("L136").Value = 8454.999899999999
$ws.Range("M136").Value = -60899.50199999999
$ws.Range("N136").Value = -13554.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2402
$ws.Range("I96").Value = 2400
$ws.Range("J96").Value = 2402.3333
$ws.Range("K96").Value = 2400
$ws.Range("L96").Value = 2402.3333
$ws.Range("M96").Value = -1027
$ws.Range("N96").Value = -5148.3333
$ws.Range("H100").Value = 33791
$ws.Range("I100").Value = 33791
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 67582
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -67041
$ws.Range("N100").ClearContents()
$ws.Range("H104").Value = 40842.25
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 40842.25
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 40842.25
$ws.Range("N104").Value = -47830.25
$ws.Range("H105").Value = 45722
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 45722
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 45722
$ws.Range("N105").Value = -52710
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("J109").Value = 40000
$ws.Range("K109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("M109").Value = -38613
$ws.Range("N109").Value = -42774
$ws.Range("H126").Value = 1670.7273
$ws.Range("I126").Value = 898.8333
$ws.Range("J126").Value = 2597
$ws.Range("K126").Value = 2696.4999
$ws.Range("L126").Value = 7791
$ws.Range("M126").Value = -226.4998999999998
$ws.Range("N126").Value = -12731
$ws.Range("H132").Value = 1747.8889
$ws.Range("I132").Value = 1074.6923
$ws.Range("J132").Value = 3498.2
$ws.Range("K132").Value = 3224.0769
$ws.Range("L132").Value = 10494.6
$ws.Range("M132").Value = -694.0769
$ws.Range("N132").Value = -15554.6
$ws.Range("H136").Value = 1477.7
$ws.Range("I136").Value = 813.0909
$ws.Range("J136").Value = 2290
$ws.Range("K136").Value = 2439.2727
$ws.Range("L136").Value = 6870
$ws.Range("M136").Value = 110.7273

